$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.147904872894287
$ws.Range("B1").Value = 1.531311273574829
$ws.Range("C1").Value = 3.561938047409058
$ws.Range("D1").Value = 3.697333574295044
$ws.Range("E1").Value = 0.9846777319908142
